# LOM3097.xlsx edit: remove the two "bare" professor-name rows (old rows 13/14,
# which only carried B/C values with no A-column label) and let everything
# below shift up by two rows. After the shift, a handful of B/C cells get
# overwritten with the (reshuffled) text that the target workbook actually
# contains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete rows 13 and 14 entirely - this shifts rows 15-25 up to 13-23,
#    which realigns every remaining row label with the row height pattern
#    used afterwards (verified against the target row-height layout).
$ws.Range("A13:C14").EntireRow.Delete() | Out-Null

# 2) Patch the B/C columns of the rows whose text content changed as part
#    of this edit (labels in column A are unaffected).
$ws.Range("B10").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C10").Value = "3577649 - Carlos Angelo Nunes"

# "01/01/2013" looks like a date to Excel's smart-entry parser, which would
# store it as a date serial (and allocate a brand-new number-format style)
# instead of reusing the existing text shared-string / cell style. Force it
# in as text with a quote-prefix, then paste-special just the *formats* from
# an untouched neighboring text cell (same column/style) back on top so the
# cell ends up with the original plain text style (s="2"/s="3"), matching
# the rest of the sheet.
$ws.Range("B13").Value = "'01/01/2013"
$ws.Range("C13").Value = "'01/01/2013"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

$ws.Range("B15").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C15").Value = "3577649 - Carlos Angelo Nunes"

$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

$ws.Range("B20").Value = "Critério`nMF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$ws.Range("C20").Value = "Critério`nMF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."

$ws.Range("B21").Value = "Não será oferecida recuperação."
$ws.Range("C21").Value = "Não será oferecida recuperação."
